$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.074.87'
$ws.Range('E2').Value = '  +0.08%  '

$ws.Range('D3').Value = '1.838.07'
$ws.Range('E3').Value = '  +0.59%  '

$ws.Range('E4').Value = '  +0.11%  '

$ws.Range('D5').Value = '244.34'
$ws.Range('E5').Value = '  +1.59%  '

$ws.Range('D6').Value = '0.6339'
$ws.Range('E6').Value = '  +2.71%  '

$ws.Range('D7').Value = '1.001'
$ws.Range('E7').Value = '  +0.10%  '

$ws.Range('D8').Value = '0.07603'
$ws.Range('E8').Value = '  +3.68%  '

$ws.Range('D9').Value = '0.2951'
$ws.Range('E9').Value = '  +1.32%  '

$ws.Range('D10').Value = '22.83'
$ws.Range('E10').Value = '  +0.85%  '

$ws.Range('D11').Value = '0.07753'
$ws.Range('E11').Value = '  +0.82%  '

$ws.Range('D12').Value = '1.834.37'
$ws.Range('E12').Value = '  +0.43%  '

$ws.Range('D13').Value = '5.004'
$ws.Range('E13').Value = '  +0.89%  '

$ws.Range('D14').Value = '0.6714'
$ws.Range('E14').Value = '  +1.58%  '

$ws.Range('D15').Value = '83.35'
$ws.Range('E15').Value = '  +2.03%  '

$ws.Range('D16').Value = '0.000009866'
$ws.Range('E16').Value = '  +9.46%  '

$ws.Range('D17').Value = '6.136'
$ws.Range('E17').Value = '  +1.88%  '

$ws.Range('D18').Value = '29.107.92'
$ws.Range('E18').Value = '  +0.23%  '

$ws.Range('B19').Value = 'Avalanche'
$ws.Range('C19').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D19').Value = '12.56'
$ws.Range('E19').Value = '  +1.67%  '

$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').Value = '227.18'
$ws.Range('E20').Value = '  +0.99%  '

$ws.Range('E21').Value = '  +0.01%  '

$ws.Range('D22').Value = '7.259'
$ws.Range('E22').Value = '  +2.17%  '

$ws.Range('D23').Value = '1.002'

$ws.Range('E24').Value = '  +0.35%  '

$ws.Range('D25').Value = '0.1410'
$ws.Range('E25').Value = '  +4.45%  '

$ws.Range('D26').Value = '8.556'
$ws.Range('E26').Value = '  +1.67%  '

$ws.Range('D27').Value = '17.98'
$ws.Range('E27').Value = '  +1.19%  '

$ws.Range('D28').Value = '1.504'
$ws.Range('E28').Value = '  +0.72%  '

$ws.Range('D29').Value = '4.128'
$ws.Range('E29').Value = '  +2.03%  '

$ws.Range('D30').Value = '4.048'
$ws.Range('E30').Value = '  +0.30%  '

$ws.Range('D31').Value = '1.204'
$ws.Range('E31').Value = '  +0.46%  '

$ws.Range('D32').Value = '0.05390'
$ws.Range('E32').Value = '  +2.45%  '

$ws.Range('D33').Value = '1.866'
$ws.Range('E33').Value = '  +1.53%  '

$ws.Range('D34').Value = '0.7513'
$ws.Range('E34').Value = '  +2.83%  '

$ws.Range('D35').Value = '1.145'
$ws.Range('E35').Value = '  -0.26%  '

$ws.Range('D36').Value = '2.670'
$ws.Range('E36').Value = '  +0.89%  '

$ws.Range('D37').Value = '1.251.91'
$ws.Range('E37').Value = '  -3.38%  '

$ws.Range('D38').Value = '0.01801'
$ws.Range('E38').Value = '  +1.16%  '

$ws.Range('E39').Value = '  +0.61%  '

$ws.Range('D40').Value = '6.567'
$ws.Range('E40').Value = '  +4.17%  '

$ws.Range('D41').Value = '0.9078'
$ws.Range('E41').Value = '  +0.69%  '

$ws.Range('E42').Value = '  +0.25%  '

$ws.Range('D43').Value = '102.91'
$ws.Range('E43').Value = '  +1.13%  '

$ws.Range('D44').Value = '1.982.46'
$ws.Range('E44').Value = '  +0.52%  '

$ws.Range('D45').Value = '0.00000000123'
$ws.Range('E45').Value = '  +2.90%  '

$ws.Range('D46').Value = '64.95'
$ws.Range('E46').Value = '  +1.90%  '

$ws.Range('E47').Value = '  +0.12%  '

$ws.Range('D48').Value = '0.4107'
$ws.Range('E48').Value = '  +3.76%  '

$ws.Range('D49').Value = '9.111'
$ws.Range('E49').Value = '  +3.56%  '

$ws.Range('D50').Value = '0.05811'
$ws.Range('E50').Value = '  +0.48%  '

$ws.Range('D51').Value = '6.781'
$ws.Range('E51').Value = '  +1.91%  '
